# Natmi following Dr Hou advice
# Recomputed NATMI ligand-receptor edge table (Agt -> Mas1) for sCs/FAPs clusters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Mas1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4683593333333333
$ws.Range("H2").Value = 1.405078
$ws.Range("I2").Value = 0.6051469521021553
$ws.Range("J2").Value = 0.6051469521021552
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6774196666666666
$ws.Range("N2").Value = 2.032259
$ws.Range("O2").Value = 0.1470994981014133
$ws.Range("P2").Value = 0.1470994981014133
$ws.Range("Q2").Value = 0.3172758234668889
$ws.Range("R2").Value = 2.855482411202
$ws.Range("S2").Value = 0.08901681293182706
$ws.Range("T2").Value = 0.08901681293182705

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Mas1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4683593333333333
$ws.Range("H3").Value = 1.405078
$ws.Range("I3").Value = 0.6051469521021553
$ws.Range("J3").Value = 0.6051469521021552
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.062384333333333
$ws.Range("N3").Value = 9.187152999999999
$ws.Range("O3").Value = 0.664986891572823
$ws.Range("P3").Value = 0.664986891572823
$ws.Range("Q3").Value = 1.434296284770444
$ws.Range("R3").Value = 12.908666562934
$ws.Range("S3").Value = 0.4024147906231803
$ws.Range("T3").Value = 0.4024147906231802

# Row 4: FAPs -> sCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Agt"
$ws.Range("C4").Value = "Mas1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4683593333333333
$ws.Range("H4").Value = 1.405078
$ws.Range("I4").Value = 0.6051469521021553
$ws.Range("J4").Value = 0.6051469521021552
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.865376
$ws.Range("N4").Value = 2.596128
$ws.Range("O4").Value = 0.1879136103257636
$ws.Range("P4").Value = 0.1879136103257636
$ws.Range("Q4").Value = 0.4053069264426667
$ws.Range("R4").Value = 3.647762337984001
$ws.Range("S4").Value = 0.1137153485471479
$ws.Range("T4").Value = 0.1137153485471479

# Row 5: sCs -> ECs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Agt"
$ws.Range("C5").Value = "Mas1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3056003333333333
$ws.Range("H5").Value = 0.916801
$ws.Range("I5").Value = 0.3948530478978448
$ws.Range("J5").Value = 0.3948530478978448
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6774196666666666
$ws.Range("N5").Value = 2.032259
$ws.Range("O5").Value = 0.1470994981014133
$ws.Range("P5").Value = 0.1470994981014133
$ws.Range("Q5").Value = 0.2070196759398889
$ws.Range("R5").Value = 1.863177083459
$ws.Range("S5").Value = 0.0580826851695863
$ws.Range("T5").Value = 0.0580826851695863

# Row 6: sCs -> FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Agt"
$ws.Range("C6").Value = "Mas1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3056003333333333
$ws.Range("H6").Value = 0.916801
$ws.Range("I6").Value = 0.3948530478978448
$ws.Range("J6").Value = 0.3948530478978448
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.062384333333333
$ws.Range("N6").Value = 9.187152999999999
$ws.Range("O6").Value = 0.664986891572823
$ws.Range("P6").Value = 0.664986891572823
$ws.Range("Q6").Value = 0.9358656730614442
$ws.Range("R6").Value = 8.422791057552999
$ws.Range("S6").Value = 0.2625721009496428
$ws.Range("T6").Value = 0.2625721009496428

# Row 7: sCs -> sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Agt"
$ws.Range("C7").Value = "Mas1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3056003333333333
$ws.Range("H7").Value = 0.916801
$ws.Range("I7").Value = 0.3948530478978448
$ws.Range("J7").Value = 0.3948530478978448
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.865376
$ws.Range("N7").Value = 2.596128
$ws.Range("O7").Value = 0.1879136103257636
$ws.Range("P7").Value = 0.1879136103257636
$ws.Range("Q7").Value = 0.2644591940586666
$ws.Range("R7").Value = 2.380132746528
$ws.Range("S7").Value = 0.07419826177861569
$ws.Range("T7").Value = 0.07419826177861569
